$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Start a new trailing paragraph after the existing last paragraph, then
# append the three runs (each InsertAfter call on the shape's full text
# range creates its own run at the end of that new paragraph).
[void]$tr.InsertAfter("`r")
[void]$tr.InsertAfter('This is Euscript ')
[void]$tr.InsertAfter('$\mathcal{A} \neq \EuScript{A}$')
[void]$tr.InsertAfter('.')
